$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$values = @(
    "fangen",
    "flüchten",
    "spinnen",
    "wehren",
    "landen",
    "fallen",
    "schenken",
    "gründen",
    "sorgen",
    "prügeln",
    "spielen",
    "nerven",
    "wirken",
    "melden",
    "boxen",
    "steuern",
    "wüten",
    "dienen",
    "nähen",
    "knarren",
    "wundern",
    "dürfen",
    "erben",
    "achten",
    "zünden",
    "sinken",
    "streichen",
    "altern",
    "sparen",
    "betteln",
    "hören",
    "ändern",
    "haben",
    "suchen",
    "bellen",
    "trennen",
    "stören",
    "schwingen",
    "husten",
    "bieten",
    "stammen",
    "feiern",
    "ärgern",
    "träumen",
    "schlagen",
    "scheinen",
    "ahnen",
    "sichern",
    "rasen",
    "formen",
    "buchen",
    "tauchen",
    "wagen",
    "schreiten",
    "kosten",
    "gelten",
    "feuern",
    "eignen",
    "mögen",
    "räumen",
    "klingen",
    "opfern",
    "leugnen",
    "scheiden",
    "quälen",
    "schwören",
    "freuen",
    "hauen",
    "brauchen",
    "herzen",
    "beißen",
    "trauen",
    "flehen",
    "weichen",
    "heißen",
    "tropfen",
    "läuten",
    "bilden",
    "seufzen",
    "filmen",
    "irren",
    "greifen",
    "jubeln",
    "kratzen",
    "jagen",
    "lügen",
    "scheitern",
    "schmecken",
    "fügen",
    "machen",
    "retten",
    "grüßen",
    "rufen",
    "lockern",
    "schulden",
    "locken",
    "geben",
    "decken",
    "fließen",
    "tragen",
    "enden",
    "wenden",
    "wachsen",
    "ehren",
    "töten",
    "lächeln",
    "reizen",
    "münzen",
    "pissen",
    "folgen",
    "zeugen",
    "schleppen",
    "gnaden",
    "sterben",
    "hupen",
    "platzen",
    "lenken",
    "knurren",
    "klettern",
    "sperren",
    "stoppen",
    "malen",
    "drucken",
    "bergen",
    "heulen",
    "stechen",
    "doppeln",
    "tollen",
    "werden",
    "schwächen",
    "siegen",
    "hacken",
    "runden",
    "bitten",
    "fahren",
    "helfen",
    "biegen",
    "rücken",
    "zögern",
    "kichern",
    "zielen",
    "graben",
    "zeichnen",
    "treten",
    "wehtun",
    "loben",
    "heilen",
    "liegen",
    "saufen",
    "leeren",
    "knien",
    "zeigen",
    "kümmern",
    "pflanzen",
    "starren",
    "stecken",
    "kehren",
    "treiben",
    "mühen",
    "pflegen",
    "bluten",
    "liefern",
    "küssen",
    "warnen",
    "messen",
    "bauen",
    "führen",
    "dringen",
    "pfeifen",
    "herrschen",
    "arten",
    "heben",
    "lesen",
    "wetten",
    "schlucken",
    "rauben",
    "spüren",
    "fesseln",
    "fischen",
    "fällen",
    "sprengen",
    "werfen",
    "füttern",
    "äußern",
    "stehlen",
    "schrecken",
    "drehen",
    "backen",
    "schützen",
    "mauern",
    "löschen",
    "planen"
)
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $values[$i]
}
